$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1772
$ws1.Range("F6").Value = 259

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1772
$ws4.Range("F7").Value = 259
